$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group")

# Drop the old "Mortality (Y/N)" column (column N) from the Group sheet.
$ws.Columns("N").Delete()

# Insert a new first column to hold the "Sample" field, pushing every
# other column (and their widths) one slot to the right.
$ws.Columns("A").Insert()
$ws.Range("A1").Value = "Sample"
$ws.Range("A1").Font.Bold = $true

# The "Group" tab is now the one the workbook opens on, instead of
# "Individual".
$ws.Select() | Out-Null
$ws.Range("A1").Select() | Out-Null
